$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is (cell address, literal new value). Price-column entries are
# prefixed with a literal leading apostrophe so Excel keeps them as text
# (these "prices" like 25.928.16 / 0.9995 are display strings, not numbers),
# then the quote-prefix cell formatting that introduces is cleared again so
# the cell style is left untouched.
$updates = @(
    ,@("D2", "'25.928.16")
    ,@("E2", "  -0.81%  ")
    ,@("D3", "'1.744.79")
    ,@("E3", "  -0.29%  ")
    ,@("D4", "'0.9995")
    ,@("E4", "  -0.07%  ")
    ,@("D5", "'246.98")
    ,@("E5", "  +4.27%  ")
    ,@("D6", "'0.9997")
    ,@("E6", "  -0.07%  ")
    ,@("D7", "'0.5047")
    ,@("E7", "  -4.65%  ")
    ,@("E8", "  -3.02%  ")
    ,@("D9", "'0.06182")
    ,@("E9", "  -0.08%  ")
    ,@("D10", "'1.752.94")
    ,@("E10", "  +0.18%  ")
    ,@("D11", "'0.07252")
    ,@("E11", "  +0.93%  ")
    ,@("D12", "'0.6535")
    ,@("E12", "  +1.29%  ")
    ,@("D13", "'15.17")
    ,@("E13", "  -1.82%  ")
    ,@("D14", "'4.639")
    ,@("E14", "  +0.10%  ")
    ,@("D15", "'77.70")
    ,@("E15", "  -1.16%  ")
    ,@("D16", "'1.0000")
    ,@("E16", "  +0.00%  ")
    ,@("D17", "'0.9992")
    ,@("E17", "  -0.09%  ")
    ,@("D18", "'25.953.12")
    ,@("E18", "  -0.32%  ")
    ,@("D19", "'11.83")
    ,@("E19", "  +0.46%  ")
    ,@("E20", "  +1.05%  ")
    ,@("D21", "'1.971.27")
    ,@("E21", "  -0.03%  ")
    ,@("D22", "'4.365")
    ,@("E22", "  +0.85%  ")
    ,@("D23", "'8.687")
    ,@("E23", "  -0.44%  ")
    ,@("D24", "'5.396")
    ,@("E24", "  +3.07%  ")
    ,@("D25", "'136.74")
    ,@("E25", "  -2.25%  ")
    ,@("D26", "'1.501")
    ,@("E26", "  -1.19%  ")
    ,@("D27", "'15.24")
    ,@("E27", "  -0.43%  ")
    ,@("D28", "'1.774")
    ,@("E28", "  -2.03%  ")
    ,@("D29", "'105.59")
    ,@("E29", "  +0.39%  ")
    ,@("D30", "'3.909")
    ,@("E30", "  +2.70%  ")
    ,@("D31", "'0.08231")
    ,@("E31", "  -0.91%  ")
    ,@("D32", "'3.642")
    ,@("E32", "  -0.12%  ")
    ,@("E33", "  +1.01%  ")
    ,@("E34", "  +0.12%  ")
    ,@("D35", "'0.9928")
    ,@("E35", "  -2.37%  ")
    ,@("D36", "'0.6186")
    ,@("E36", "  -2.51%  ")
    ,@("E37", "  +1.64%  ")
    ,@("D38", "'0.01607")
    ,@("E38", "  -1.04%  ")
    ,@("D39", "'1.920")
    ,@("E39", "  -2.99%  ")
    ,@("D40", "'0.9996")
    ,@("E40", "  -0.03%  ")
    ,@("D41", "'100.02")
    ,@("E41", "  -2.57%  ")
    ,@("D42", "'0.3917")
    ,@("E42", "  -0.54%  ")
    ,@("D43", "'0.7572")
    ,@("E43", "  +0.54%  ")
    ,@("D44", "'5.009")
    ,@("E44", "  -0.86%  ")
    ,@("D45", "'0.1145")
    ,@("E45", "  -0.90%  ")
    ,@("D46", "'6.294")
    ,@("E46", "  -1.24%  ")
    ,@("D47", "'55.50")
    ,@("E47", "  +1.57%  ")
    ,@("D48", "'0.05257")
    ,@("E48", "  -1.75%  ")
    ,@("D49", "'30.62")
    ,@("E49", "  -1.55%  ")
    ,@("D50", "'7.572")
    ,@("E50", "  -0.41%  ")
    ,@("D51", "'0.3430")
    ,@("E51", "  -1.71%  ")
)

foreach ($u in $updates) {
    $addr = $u[0]
    $val = $u[1]
    $cell = $ws.Range($addr)
    $cell.Value = $val
    $cell.Style = "Normal"
}
